# Atualização automática de ESTRELA.xlsx
$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the sheet that is no longer needed
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Keep the original sheet (now "PAINEIS DARQ") as the active one, since
# deleting a sheet can shift activation to a neighboring tab.
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
